$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Widen column E for the new note text
$ws1.Columns.Item(5).ColumnWidth = 15.42578125

# Add the entrance instructions next to the relevant lectures
$ws1.Range("E10").Value = "використовуйте  основний вхід корпусу А"
$ws1.Range("E10").WrapText = $true
$ws1.Range("E10").HorizontalAlignment = -4108
$ws1.Range("E10").VerticalAlignment = -4108

$ws1.Range("E22").Value = "використовуйте  додатковий вхід корпусу П"
$ws1.Range("E22").WrapText = $true
$ws1.Range("E22").HorizontalAlignment = -4108
$ws1.Range("E22").VerticalAlignment = -4108

# Increase row 22 height to fit the longer note
$ws1.Rows.Item(22).RowHeight = 45.75

# Select C4 as the active cell like the saved view
$ws1.Range("C4").Select()

# Match the print scale change
$ws1.PageSetup.Zoom = 80
